$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0" and J1 = "IF", matching the style
#     (bold, centered, bordered) already used by the other header cells. ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J1").Value = "IF"

$ws.Application.CutCopyMode = $false

# --- Data rows 2..34: I = constant 1, J = copy of H's value. ---
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hVal
}
